$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Expiry date text for row 11 (Kotireddy) first, so shared-string order matches source (index 50)
$ws.Range("F11").Value = "25th JAN 2026"

# New policy number text for Alekya Kotireddy term insurance (row 10) - shared string index 51
$ws.Range("C10").Value = "K3220367"

# Payment till / Cover values for row 10 (Alekya) - copy date-style number format from existing similarly formatted cell
$ws.Range("J10").Value = 46009
$ws.Range("K10").Value = 65368
$ws.Range("J10:K10").NumberFormat = $ws.Range("K9").NumberFormat
$ws.Range("J10:K10").HorizontalAlignment = $ws.Range("K9").HorizontalAlignment

# Row 11 (Kotireddy) - policy number, payment till / cover values
$ws.Range("C11").Value = 170557110
$ws.Range("J11").Value = 21210
$ws.Range("K11").Value = 58101
$ws.Range("J11:K11").NumberFormat = $ws.Range("K9").NumberFormat
$ws.Range("J11:K11").HorizontalAlignment = $ws.Range("K9").HorizontalAlignment

# Update the view: scroll so column C is the left-most visible column and select J10
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("J10").Select()
